$d = $word.ActiveDocument

# Locate the "url" meta paragraph precisely (search restricted to a single
# paragraph's Range so we do not hit the other "content=" occurrences).
$targetPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like '*name="url"*') {
        $targetPara = $p
    }
}

$pr = $targetPara.Range.Duplicate
$pr.Find.Execute('" content="', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterLabel = $pr.End

$pr2 = $targetPara.Range.Duplicate
$pr2.Find.Execute('"' + [char]62, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$closeStart = $pr2.Start
$closeEnd = $pr2.End

# Step 1: rewrite the whole URL content in one shot (this is expected to
# collapse runs 2-5 into a single run covering [afterLabel, closeEnd)).
$whole = $d.Range($afterLabel, $closeEnd)
$whole.Text = ' https://starsindust.github.io/ Enlightenment/Articles/2025/4_Game_Maker_2/3_Creating_Objects/3_Creating_Objects.html ">'

Write-Output "after step1"
Write-Output $targetPara.Range.Text
